# Update Metadata sheet: Last Updated timestamp
$wb = $excel.ActiveWorkbook
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("A2").Value = "05 Nov 2025, 03:28 PM"

# Update Stock List sheet: rows shift up by one (oldest entry dropped),
# a new entry (TRAVELFOOD) is appended at the bottom (row 76).
$wsStock = $wb.Worksheets.Item("Stock List")

$stockData = New-Object "object[]" 75
$stockData[0] = @("NIFTYCASE", "NIFTYCASE", 10.19, -0.5854, 0.0)
$stockData[1] = @("MOMENTUM30", "MOMENTUM30", 31.54, -0.6614, 0.0)
$stockData[2] = @("CANHLIFE", "CANHLIFE", 118.46, 0.6286, 11253.7)
$stockData[3] = @("FLEXIADD", "FLEXIADD", 10.64, -1.0233, 0.0)
$stockData[4] = @("MOENERGY", "MOENERGY", 36.3, -0.6568, 0.0)
$stockData[5] = @("MONIFTY100", "MONIFTY100", 26.49, 0.3409, 0.0)
$stockData[6] = @("RUBICON", "RUBICON", 652.65, -0.1453, 10752.4289)
$stockData[7] = @("CRAMC", "CRAMC", 317.2, 2.3226, 6325.5208)
$stockData[8] = @("LGEINDIA", "LGEINDIA", 1633.4, -0.946, 110870.6825)
$stockData[9] = @("TATACAP", "TATACAP", 329.3, 0.1521, 139783.5374)
$stockData[10] = @("ELIQUID", "ELIQUID", 1004.85, 0.0408, 0.0)
$stockData[11] = @("WEWORK", "WEWORK", 632.15, -2.4008, 8472.2803)
$stockData[12] = @("GROWWRLTY", "GROWWRLTY", 10.8, -0.4608, 0.0)
$stockData[13] = @("ADVANCE", "ADVANCE", 130.05, -5.2666, 836.0358)
$stockData[14] = @("OMFREIGHT", "OMFREIGHT", 88.9, -0.5926, 299.3747)
$stockData[15] = @("GLOTTIS", "GLOTTIS", 72.74, -0.8587, 672.1394)
$stockData[16] = @("FABTECH", "FABTECH", 237.72, 0.4734, 1056.6843)
$stockData[17] = @("PACEDIGITK", "PACEDIGITK", 218.85, 0.1327, 4723.9063)
$stockData[18] = @("JAINREC", "JAINREC", 377.25, 1.2208, 13018.3623)
$stockData[19] = @("EPACKPEB", "EPACKPEB", 301.45, 1.979, 3028.1254)
$stockData[20] = @("BMWVENTLTD", "BMWVENTLTD", 69.25, 0.0, 600.5014)
$stockData[21] = @("STYL", "STYL", 372.4, -0.8388, 6025.649)
$stockData[22] = @("JARO", "JARO", 621.5, -1.4821, 1377.0134)
$stockData[23] = @("SOLARWORLD", "SOLARWORLD", 309.1, -0.6269, 2679.0517)
$stockData[24] = @("ARSSBL", "ARSSBL", 537.3, 4.7266, 3370.2277)
$stockData[25] = @("GANESHCP", "GANESHCP", 274.4, -2.7984, 1108.9312)
$stockData[26] = @("ATLANTAELE", "ATLANTAELE", 1003.05, -1.7436, 7713.116)
$stockData[27] = @("GKENERGY", "GKENERGY", 213.85, -0.7933, 4337.2472)
$stockData[28] = @("SAATVIKGL", "SAATVIKGL", 528.2, -1.3079, 6713.6863)
$stockData[29] = @("IVALUE", "IVALUE", 281.45, -0.3364, 1506.8799)
$stockData[30] = @("VMSTMT", "VMSTMT", 70.03, -0.9056, 347.5674)
$stockData[31] = @("EUROPRATIK", "EUROPRATIK", 321.75, 0.8147, 3288.285)
$stockData[32] = @("SHRINGARMS", "SHRINGARMS", 229.31, -1.2616, 2211.284)
$stockData[33] = @("DEVX", "DEVX", 44.53, -0.3803, 401.605)
$stockData[34] = @("URBANCO", "URBANCO", 148.9, -2.0459, 21380.5798)
$stockData[35] = @("SML100CASE", "SML100CASE", 10.36, -0.7663, 0.0)
$stockData[36] = @("AONEGOLD", "AONEGOLD", 11.28, -0.2653, 0.0)
$stockData[37] = @("ELM250", "ELM250", 16.72, 0.1797, 0.0)
$stockData[38] = @("AMANTA", "AMANTA", 122.52, 1.407, 475.7372)
$stockData[39] = @("CPEDU", "CPEDU", 315.9, 1.8539, 574.7149)
$stockData[40] = @("AHCL", "AHCL", 139.27, 3.1706, 740.2409)
$stockData[41] = @("STLNETWORK", "STLNETWORK", 26.59, -0.412, 1297.3822)
$stockData[42] = @("VIKRAN", "VIKRAN", 98.05, -1.783, 2528.8166)
$stockData[43] = @("MANUFGBEES", "MANUFGBEES", 151.77, -1.011, 0.0)
$stockData[44] = @("MEIL", "MEIL", 461.15, -0.7319, 1274.1632)
$stockData[45] = @("GROWWNXT50", "GROWWNXT50", 70.29, -0.4109, 0.0)
$stockData[46] = @("SHREEJISPG", "SHREEJISPG", 270.05, -0.7899, 4399.6074)
$stockData[47] = @("GEMAROMA", "GEMAROMA", 219.52, -0.876, 1146.7097)
$stockData[48] = @("PATELRMART", "PATELRMART", 219.31, -1.0646, 732.507)
$stockData[49] = @("VIKRAMSOLR", "VIKRAMSOLR", 322.0, -1.5892, 11647.2884)
$stockData[50] = @("LTGILTCASE", "LTGILTCASE", 29.67, 0.2365, 0.0)
$stockData[51] = @("REGAAL", "REGAAL", 89.13, -0.8675, 915.5742)
$stockData[52] = @("BLUESTONE", "BLUESTONE", 711.95, 0.1266, 10773.2539)
$stockData[53] = @("MOSILVER", "MOSILVER", 145.9, -1.5054, 0.0)
$stockData[54] = @("ALLTIME", "ALLTIME", 308.75, 2.66, 2022.5526)
$stockData[55] = @("JSWCEMENT", "JSWCEMENT", 134.98, -0.4793, 18402.6999)
$stockData[56] = @("SBILIQETF", "SBILIQETF", 1012.94, 0.0296, 0.0)
$stockData[57] = @("HILINFRA", "HILINFRA", 77.23, -0.3998, 0.0)
$stockData[58] = @("GROWWPOWER", "GROWWPOWER", 10.28, -0.9634, 0.0)
$stockData[59] = @("LOTUSDEV", "LOTUSDEV", 177.82, 0.3669, 8690.485)
$stockData[60] = @("MBEL", "MBEL", 450.2, -0.7714, 2572.8126)
$stockData[61] = @("LAXMIINDIA", "LAXMIINDIA", 145.62, -1.1942, 761.1248)
$stockData[62] = @("CPPLUS", "CPPLUS", 1322.1, -0.264, 15497.9053)
$stockData[63] = @("SHANTIGOLD", "SHANTIGOLD", 241.57, -1.6409, 1741.6231)
$stockData[64] = @("MOGOLD", "MOGOLD", 119.65, -0.5403, 0.0)
$stockData[65] = @("BRIGHOTEL", "BRIGHOTEL", 82.39, -0.9855, 3129.5229)
$stockData[66] = @("INDIQUBE", "INDIQUBE", 212.64, -0.7561, 4465.6847)
$stockData[67] = @("EBGNG", "EBGNG", 346.65, 3.2311, 3952.2092)
$stockData[68] = @("LIQGRWBEES", "LIQGRWBEES", 1014.74, 0.0246, 0.0)
$stockData[69] = @("CHEMBONDCH", "CHEMBONDCH", 153.35, -1.6987, 412.459)
$stockData[70] = @("GROWWNIFTY", "GROWWNIFTY", 10.29, -0.3872, 0.0)
$stockData[71] = @("ANTHEM", "ANTHEM", 702.25, -0.1209, 39439.0658)
$stockData[72] = @("QUALITY30", "QUALITY30", 21.05, -0.8945, 0.0)
$stockData[73] = @("SMARTWORKS", "SMARTWORKS", 606.65, 2.0867, 6931.2448)
$stockData[74] = @("TRAVELFOOD", "TRAVELFOOD", 1316.3, 0.1141, 17332.9705)

for ($i = 0; $i -lt 75; $i++) {
    $row = $i + 2
    $rec = $stockData[$i]
    $wsStock.Cells.Item($row, 2).Value = $rec[0]
    $wsStock.Cells.Item($row, 3).Value = $rec[1]
    $wsStock.Cells.Item($row, 4).Value = $rec[2]
    $wsStock.Cells.Item($row, 5).Value = $rec[3]
    $wsStock.Cells.Item($row, 8).Value = $rec[4]
}
